$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 14 (leo) values
$ws.Range("B14").Value = 78
$ws.Range("C14").Value = 722

# Add new row 15 for alexandra (A15, B15, C15); D15 already has policy text
$ws.Range("A15").Value = "alexandra"
$ws.Range("B15").Value = 29
$ws.Range("C15").Value = 486
